# Fruta / hortaliza, semanal
# Insert 3 new weekly price rows for "Clemenuless" (mandarina) above the
# existing data block, shifting the remaining rows (525-612) down to
# 528-615.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert three blank rows starting at row 525; this pushes the former
# rows 525-612 down to 528-615 and extends the used range accordingly.
$ws.Rows("525:527").Insert()

# Columns that are constant for every data row in this single-vendor,
# single-product sheet.
$mercadoId = 2
$mercado   = "Comercializadora del Agro de Limarí"
$region    = "Coquimbo"
$codreg    = 4
$tipo      = "Fruta"
$productoId = 100102
$producto  = "Cítricos"
$categoriaId = 100102004
$categoria = "Mandarina"
$origen    = "Provincia de Limarí"

function Set-PrecioRow($row, $fecha, $variedad, $calidad, $volumen, $precioMin, $precioMax, $precioProm, $unidad, $precioKg, $kgUnidad) {
    $ws.Cells.Item($row, 1).Value  = $mercadoId
    $ws.Cells.Item($row, 2).Value  = $mercado
    $ws.Cells.Item($row, 3).Value  = $region
    $ws.Cells.Item($row, 4).Value  = $fecha
    $ws.Cells.Item($row, 5).Value  = $codreg
    $ws.Cells.Item($row, 6).Value  = $tipo
    $ws.Cells.Item($row, 7).Value  = $productoId
    $ws.Cells.Item($row, 8).Value  = $producto
    $ws.Cells.Item($row, 9).Value  = $categoriaId
    $ws.Cells.Item($row, 10).Value = $categoria
    $ws.Cells.Item($row, 11).Value = $variedad
    $ws.Cells.Item($row, 12).Value = $calidad
    $ws.Cells.Item($row, 13).Value = $volumen
    $ws.Cells.Item($row, 14).Value = $precioMin
    $ws.Cells.Item($row, 15).Value = $precioMax
    $ws.Cells.Item($row, 16).Value = $precioProm
    $ws.Cells.Item($row, 17).Value = $unidad
    $ws.Cells.Item($row, 18).Value = $origen
    $ws.Cells.Item($row, 19).Value = $precioKg
    $ws.Cells.Item($row, 20).Value = $kgUnidad
}

# New rows, all dated 2023-07-20 (serial 45127).
Set-PrecioRow 525 "7/20/2023" "Clemenuless" "Especial" 1100 6000 6500 6250 "$/bandeja 10 kilos" 625 10
Set-PrecioRow 526 "7/20/2023" "Clemenuless" "Primera"  1200 4000 4500 4250 "$/bandeja 10 kilos" 425 10
Set-PrecioRow 527 "7/20/2023" "Clemenuless" "Segunda"  1000 2000 2500 2250 "$/bandeja 10 kilos" 225 10
